# Updates the SQL queries on the active worksheet so that the JOIN
# conditions use the fully-qualified study_id / participant_id columns
# instead of the generic "id" column, matching the updated C3DC schema.
#
# Old:
#   df_participant prt ON std.id = prt."study.id"
#   df_diagnoses dgn ON prt.id = dgn."participant.id"
#   df_treatments trt ON prt.id = trt."participant.id"
#   df_treatment_resp trr ON prt.id = trr."participant.id"
#   df_survival srv ON prt.id = srv."participant.id"
#   df_reference_files rfs ON std.id = rfs."study.id"
#
# New:
#   df_participant prt ON std.study_id = prt."study.study_id"
#   df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"
#   df_treatments trt ON prt.participant_id = trt."participant.participant_id"
#   df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"
#   df_survival srv ON prt.participant_id = srv."participant.participant_id"
#   df_reference_files rfs ON std.study_id = rfs."study.study_id"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (old, new) substring replacements applied to every
# cell that contains SQL text on the sheet.
$replacements = @(
    @{ Old = 'df_participant prt ON std.id = prt."study.id"'; New = 'df_participant prt ON std.study_id = prt."study.study_id"' },
    @{ Old = 'df_diagnoses dgn ON prt.id = dgn."participant.id"'; New = 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"' },
    @{ Old = 'df_treatments trt ON prt.id = trt."participant.id"'; New = 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"' },
    @{ Old = 'df_treatment_resp trr ON prt.id = trr."participant.id"'; New = 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"' },
    @{ Old = 'df_survival srv ON prt.id = srv."participant.id"'; New = 'df_survival srv ON prt.participant_id = srv."participant.participant_id"' },
    @{ Old = 'df_reference_files rfs ON std.id = rfs."study.id"'; New = 'df_reference_files rfs ON std.study_id = rfs."study.study_id"' }
)

# Update the query cells in the same relative order the queries appear
# in the workbook's shared string table (C2, B2, B3, B4, B5, B6, B7) so
# that the underlying shared-string ordering is preserved on save.
$cellAddresses = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellAddresses) {
    $cell = $ws.Range($addr)
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Contains('LEFT JOIN')) {
        $newVal = $val
        foreach ($rep in $replacements) {
            $newVal = $newVal.Replace($rep.Old, $rep.New)
        }
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
